$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions): row 2 is a full replacement - the old
#     event ("杭州·第三届日夜国乙only") was swapped out for a brand new one
#     ("杭州·Aniidol偶像剧场·端午甜咸对决"). Every field in the row changes. ---
$ws1 = $wb.Worksheets.Item("展览")

# B2 holds a plain text date like "2024-06-09"; Excel will happily
# auto-coerce a "YYYY-MM-DD"-shaped string into a real date serial,
# which the source file does NOT want (it is stored as literal text).
# Force the cell to Text, write the literal string, then strip the
# temporary number format again so no stray style sticks around.
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "2024-06-10"
$ws1.Range("B2").ClearFormats()

$ws1.Range("C2").Value = "杭州·Aniidol偶像剧场·端午甜咸对决"
$ws1.Range("D2").Value = "同协路288号 1928创意园"
$ws1.Range("E2").Value = "2024.06.10 12:00-06.10 20:00"
$ws1.Range("F2").Value = 33
$ws1.Range("G2").Value = 88
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=86665"
$ws1.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202406/rGOmIh0e1717400362659.jpeg"

# --- Sheet "展览": bump the remaining "想去人数" (interested-count) figures ---
$ws1.Range("F3").Value = 806
$ws1.Range("F4").Value = 2457
$ws1.Range("F6").Value = 443
$ws1.Range("F7").Value = 267
$ws1.Range("F8").Value = 163
$ws1.Range("F9").Value = 410
$ws1.Range("F10").Value = 1143
$ws1.Range("F11").Value = 502
$ws1.Range("F12").Value = 266
$ws1.Range("F13").Value = 104
$ws1.Range("F14").Value = 317
$ws1.Range("F15").Value = 5205
$ws1.Range("F17").Value = 1455
$ws1.Range("F18").Value = 3812
$ws1.Range("F19").Value = 368
$ws1.Range("F20").Value = 235
$ws1.Range("F21").Value = 308
$ws1.Range("F22").Value = 4227
$ws1.Range("F23").Value = 5623
$ws1.Range("F25").Value = 1000
$ws1.Range("F26").Value = 605
$ws1.Range("F27").Value = 3531
$ws1.Range("F28").Value = 429
$ws1.Range("F30").Value = 163
$ws1.Range("F31").Value = 105
$ws1.Range("F32").Value = 936
$ws1.Range("F33").Value = 1278
$ws1.Range("F34").Value = 94
$ws1.Range("F35").Value = 130
$ws1.Range("F36").Value = 1506
$ws1.Range("F37").Value = 168
$ws1.Range("F38").Value = 1532
$ws1.Range("F39").Value = 100
$ws1.Range("F40").Value = 989
$ws1.Range("F41").Value = 1039
$ws1.Range("F42").Value = 565
$ws1.Range("F43").Value = 72
$ws1.Range("F44").Value = 117
$ws1.Range("F45").Value = 2617
$ws1.Range("F46").Value = 107
$ws1.Range("F47").Value = 218
$ws1.Range("F48").Value = 387
$ws1.Range("F49").Value = 3809

# --- Sheet "演出": bump "想去人数" (interested-count) figures ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 1108
$ws2.Range("F6").Value = 38
$ws2.Range("F16").Value = 3
$ws2.Range("F18").Value = 5
$ws2.Range("F21").Value = 44
$ws2.Range("F22").Value = 57
$ws2.Range("F23").Value = 21

# --- Sheet "本地生活": bump "想去人数" (interested-count) figures ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3376

# --- Sheet "全部类型": bump "想去人数" (interested-count) figures ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3376
$ws4.Range("F3").Value = 806
$ws4.Range("F4").Value = 2457
$ws4.Range("F6").Value = 443
$ws4.Range("F7").Value = 267
$ws4.Range("F8").Value = 1108
$ws4.Range("F9").Value = 163
$ws4.Range("F10").Value = 410
$ws4.Range("F11").Value = 1143
$ws4.Range("F12").Value = 502
$ws4.Range("F13").Value = 266
$ws4.Range("F14").Value = 104
$ws4.Range("F15").Value = 317
$ws4.Range("F16").Value = 5205
$ws4.Range("F18").Value = 1455
$ws4.Range("F19").Value = 4228
$ws4.Range("F20").Value = 5623
$ws4.Range("F22").Value = 1000
$ws4.Range("F23").Value = 605
$ws4.Range("F24").Value = 3531
$ws4.Range("F25").Value = 429
$ws4.Range("F27").Value = 163
$ws4.Range("F28").Value = 105
$ws4.Range("F29").Value = 936
$ws4.Range("F30").Value = 1278
$ws4.Range("F31").Value = 94
$ws4.Range("F32").Value = 130
$ws4.Range("F33").Value = 1506
$ws4.Range("F34").Value = 168
$ws4.Range("F35").Value = 1532
$ws4.Range("F36").Value = 3
$ws4.Range("F37").Value = 989
$ws4.Range("F38").Value = 5
$ws4.Range("F39").Value = 565
$ws4.Range("F41").Value = 72
$ws4.Range("F42").Value = 57
$ws4.Range("F43").Value = 2618
$ws4.Range("F44").Value = 21
$ws4.Range("F45").Value = 107
$ws4.Range("F46").Value = 218
$ws4.Range("F47").Value = 387
$ws4.Range("F49").Value = 3809
